$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name/title
$ws.Name = "Through 2022-02-19"

# Update header label for February
$ws.Range("A3").Value = "February (through 02-19)"

# Update January I2 value
$ws.Range("I2").Value = 159

# Update February row (row 3)
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 25
$ws.Range("D3").Value = 41
$ws.Range("E3").Value = 39
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = 47
$ws.Range("H3").Value = 85
$ws.Range("I3").Value = 92

# Update Total row (row 4)
$ws.Range("B4").Value = 33
$ws.Range("C4").Value = 76
$ws.Range("D4").Value = 116
$ws.Range("E4").Value = 125
$ws.Range("F4").Value = 69
$ws.Range("G4").Value = 121
$ws.Range("H4").Value = 302
$ws.Range("I4").Value = 251
